$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the active selection from T3 to T2 (matches updated <selection activeCell="T2" sqref="T2"/>)
$ws.Range("T2").Select() | Out-Null

# Update the value in T2 from 284641 to 287887
$ws.Range("T2").Value = 287887
